$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.911.39"
$ws.Range("E2").Value = "  +3.57%  "
$ws.Range("D3").Value = "'3.717.61"
$ws.Range("E3").Value = "  +8.26%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'589.37"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'181.14"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").Value = "'3.710.77"
$ws.Range("E7").Value = "  +8.28%  "
$ws.Range("E8").Value = "  +3.88%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("E11").Value = "  +4.44%  "
$ws.Range("E12").Value = "  +3.03%  "
$ws.Range("D13").Value = "'0.0000289"
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").Value = "'4.317.34"
$ws.Range("E14").Value = "  +8.30%  "
$ws.Range("D15").Value = "'682.89"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "'9.05"
$ws.Range("E16").Value = "  +4.09%  "
$ws.Range("D17").Value = "'3.708.71"
$ws.Range("E17").Value = "  +7.96%  "
$ws.Range("D18").Value = "'71.976.64"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("D20").Value = "'18.17"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").Value = "'11.71"
$ws.Range("E21").Value = "  +3.00%  "
$ws.Range("D22").Value = "'6.41"
$ws.Range("E22").Value = "  +18.99%  "
$ws.Range("E23").Value = "  +3.56%  "
$ws.Range("D24").Value = "'17.85"
$ws.Range("E24").Value = "  +4.64%  "
$ws.Range("D25").Value = "'103.92"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("E26").Value = "  +3.51%  "
$ws.Range("D27").Value = "'2.85"
$ws.Range("E27").Value = "  +5.15%  "
$ws.Range("D28").Value = "'10.38"
$ws.Range("E28").Value = "  +6.65%  "
$ws.Range("E29").Value = "  +5.79%  "
$ws.Range("D30").Value = "'9.34"
$ws.Range("E30").Value = "  +5.93%  "
$ws.Range("E31").Value = "  +6.68%  "
$ws.Range("D32").Value = "'4.21"
$ws.Range("E32").Value = "  +11.30%  "
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("E34").Value = "  +3.74%  "
$ws.Range("D35").Value = "'564.46"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").Value = "'59.74"
$ws.Range("E36").Value = "  +2.57%  "
$ws.Range("D37").Value = "'3.773.37"
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").Value = "'0.145"
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("E40").Value = "  +4.99%  "
$ws.Range("D41").Value = "'35.90"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("D42").Value = "'3.48"
$ws.Range("E42").Value = "  +5.15%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "'2.82"
$ws.Range("E43").Value = "  +4.18%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0465"
$ws.Range("E44").Value = "  +9.30%  "
$ws.Range("D45").Value = "'0.354"
$ws.Range("E45").Value = "  +4.99%  "
$ws.Range("E46").Value = "  +8.37%  "
$ws.Range("D47").Value = "'3.40"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").Value = "'135.92"
$ws.Range("E51").Value = "  +3.57%  "
